$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamps for the remaining first 8 rows (Sequence 1-8)
$ws.Range("A2").Value = "2025-11-27T13:07:24.620Z"
$ws.Range("A3").Value = "2025-11-27T13:07:24.662Z"
$ws.Range("A4").Value = "2025-11-27T13:07:24.670Z"
$ws.Range("A5").Value = "2025-11-27T13:07:24.674Z"
$ws.Range("A6").Value = "2025-11-27T13:07:24.679Z"
$ws.Range("A7").Value = "2025-11-27T13:09:05.762Z"
$ws.Range("A8").Value = "2025-11-27T13:10:42.027Z"
$ws.Range("A9").Value = "2025-11-27T13:11:26.371Z"

# Remove the last two rows (Sequence 9 and Sequence 10 alerts)
$ws.Rows("10:11").Delete()
